# TC_51987 - Update slot card related test case / test data on the
# "Add Devices" sheet: the 4th/18th "Other Slot Cards" occupied count
# label moves from N10/N11, and the two related Slot-Card-occupied
# flags (I10/I11) flip from 1 (occupied) to 0 (free).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")
$ws.Activate()

# Slot card position no longer marked as occupied.
$ws.Range("I10").Value = 0
$ws.Range("I11").Value = 0

# "Other Slot Cards" counter text updated from "(4 of 6)" to "(4 of 18)".
$ws.Range("N10").Value = "Other Slot Cards  (4 of 18)"
$ws.Range("N11").Value = "Other Slot Cards  (4 of 18)"

# Leave the selection on the cell that was last edited, matching the
# recorded interaction.
$ws.Range("N10").Select()
